# Update symbol list (crypto prices and a couple of swapped/updated rows)
# as described by the commit "Updated symbol list on Sat Dec 17 11:27:43 UTC 2022
# with GitHub Actions".
#
# The Price column (D) stores numeric-looking values as TEXT (they were
# originally inline strings, e.g. "238.56"). To keep them stored as text
# rather than letting Excel auto-convert them to numbers, each Price cell
# is temporarily switched to the "@" (Text) number format before the value
# is assigned, and the style is reset back to Normal afterwards so no stray
# cell formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Simple price updates (column D) ---
Set-TextValue "D2"  "238.37"
Set-TextValue "D3"  "21.72"
Set-TextValue "D4"  "5.456"
Set-TextValue "D6"  "6.490"
Set-TextValue "D7"  "3.352"
Set-TextValue "D8"  "1.078"
Set-TextValue "D9"  "0.7912"
Set-TextValue "D11" "0.07329"
Set-TextValue "D12" "0.03201"
Set-TextValue "D13" "0.02965"
Set-TextValue "D14" "0.09258"
Set-TextValue "D15" "0.001660"
Set-TextValue "D16" "3.261"
Set-TextValue "D17" "0.04781"
Set-TextValue "D18" "0.0005743"
Set-TextValue "D19" "0.006231"
Set-TextValue "D20" "0.005114"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.904"
Set-TextValue "D26" "0.1055"
Set-TextValue "D27" "0.0004012"
Set-TextValue "D40" "0.04145"
Set-TextValue "D41" "0.006961"

# --- Row 42 / Row 43 content swap (CEJI <-> BKEXToken) plus new values ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1043"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003011"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price / label updates ---
Set-TextValue "D44" "0.009884"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"

Set-TextValue "D45" "0.00005441"
Set-TextValue "D47" "0.6756"
Set-TextValue "D48" "0.03736"
Set-TextValue "D49" "0.00002101"

Write-Host "Applied symbol list update."
